$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column values are kept as text (avoids Excel auto-numeric parsing)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "68.249.37"
$ws.Range("E2").Value = "  +1.83%  "
$ws.Range("D3").Value = "3.929.10"
$ws.Range("E3").Value = "  +0.96%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "487.01"
$ws.Range("E5").Value = "  +4.46%  "
$ws.Range("D6").Value = "146.74"
$ws.Range("E6").Value = "  +1.94%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").Value = "0.732"
$ws.Range("E9").Value = "  -0.61%  "
$ws.Range("E10").Value = "  +3.64%  "
$ws.Range("D11").Value = "0.0000360"
$ws.Range("E11").Value = "  +7.77%  "
$ws.Range("D12").Value = "42.92"
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").Value = "10.70"
$ws.Range("E13").Value = "  +3.27%  "
$ws.Range("D14").Value = "4.552.64"
$ws.Range("E14").Value = "  +0.77%  "
$ws.Range("D15").Value = "14.91"
$ws.Range("E15").Value = "  -1.76%  "
$ws.Range("D16").Value = "3.929.57"
$ws.Range("E16").Value = "  +0.87%  "
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").Value = "20.14"
$ws.Range("E18").Value = "  +1.18%  "
$ws.Range("E19").Value = "  -1.61%  "
$ws.Range("D20").Value = "68.349.84"
$ws.Range("E20").Value = "  +1.64%  "
$ws.Range("D21").Value = "447.93"
$ws.Range("E21").Value = "  +4.04%  "
$ws.Range("E22").Value = "  +0.65%  "
$ws.Range("D23").Value = "3.41"
$ws.Range("E23").Value = "  +1.81%  "
$ws.Range("D24").Value = "88.69"
$ws.Range("E24").Value = "  +0.50%  "
$ws.Range("D25").Value = "11.65"
$ws.Range("E25").Value = "  +15.92%  "
$ws.Range("D26").Value = "10.99"
$ws.Range("E26").Value = "  +14.68%  "
$ws.Range("E27").Value = "  +2.75%  "
$ws.Range("D28").Value = "39.08"
$ws.Range("E28").Value = "  +0.95%  "
$ws.Range("E29").Value = "  +2.28%  "
$ws.Range("D30").Value = "13.53"
$ws.Range("E30").Value = "  -1.26%  "
$ws.Range("E31").Value = "  +0.75%  "
$ws.Range("D32").Value = "692.43"
$ws.Range("E32").Value = "  -5.60%  "
$ws.Range("E33").Value = "  +5.16%  "
$ws.Range("D34").Value = "0.0₃0947"
$ws.Range("E34").Value = "  +21.20%  "
$ws.Range("D35").Value = "41.78"
$ws.Range("E35").Value = "  -2.67%  "
$ws.Range("D36").Value = "58.99"
$ws.Range("E36").Value = "  +1.84%  "
$ws.Range("D37").Value = "5.80"
$ws.Range("E37").Value = "  +7.68%  "
$ws.Range("E38").Value = "  -4.37%  "
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("E40").Value = "  +0.78%  "
$ws.Range("D41").Value = "0.373"
$ws.Range("E41").Value = "  +10.75%  "
$ws.Range("E42").Value = "  +12.99%  "
$ws.Range("D43").Value = "3.08"
$ws.Range("E43").Value = "  -4.29%  "
$ws.Range("D44").Value = "2.95"
$ws.Range("E44").Value = "  +5.72%  "
$ws.Range("E45").Value = "  +1.99%  "
$ws.Range("E46").Value = "  -0.16%  "
$ws.Range("E47").Value = "  +1.24%  "
$ws.Range("E48").Value = "  -0.53%  "
$ws.Range("D49").Value = "146.09"
$ws.Range("E49").Value = "  +1.87%  "
$ws.Range("D50").Value = "3.14"
$ws.Range("E50").Value = "  -0.26%  "
$ws.Range("E51").Value = "  -1.56%  "
